# CII.xlsx refresh - July 5th data update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new first column ("Name") - shifts existing Status/Period/Indicator/Pathway/Description
#    from A:E to B:F, carrying their widths/styles/merges along automatically.
$ws.Columns("A:A").Insert()
$ws.Columns("A:A").ColumnWidth = 14.83

# 2. Populate the two new rows (6 and 7) by copying the formatting of row 5 (a data row)
#    across the used range, so the new rows get the same visual style (s=4) without creating
#    brand-new style table entries.
$ws.Range("A5:F5").Copy()
$ws.Range("A6:F6").PasteSpecial(-4122)
$ws.Range("A7:F7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3. Header row (row 2) - add "Name" header in new column A
$ws.Range("A2").Value = "Name"
$ws.Range("B2").Value = "Status"
$ws.Range("C2").Value = "Period"
$ws.Range("D2").Value = "Indicator"
$ws.Range("E2").Value = "Pathway"
$ws.Range("F2").Value = "Description"

# 4. Data rows
# Row 4
$ws.Range("A4").Value = "CII platform"
$ws.Range("B4").Value = "In Process"
$ws.Range("C4").Value = "2021-First Half"
$ws.Range("D4").Value = "Tools and platforms developed"
$ws.Range("E4").Value = "Reduced barriers in auto sector"
$ws.Range("F4").Value = "First meeting of the CII platform was held on March 1, 2021"

# Row 5
$ws.Range("A5").Value = "Bilateral consultations with industry members"
$ws.Range("B5").Value = "Completed"
$ws.Range("C5").Value = "2021-Second Half"
$ws.Range("D5").Value = "Convenings/workshops"
$ws.Range("E5").Value = "Unlocking policy and regulatory barriers"
$ws.Range("F5").Value = "Between April and July 2021, CII undertook 23 bilateral consultations with industry members to identify barriers faced by the different players across EV value chain and potential solutions required to address these challenges.  Some of the topics these consultations discussed are: emerging consensus on focus areas to accelerate EV adoption; Barriers, priorities, potential way forward for the charging infrastructure providers and Discoms and OEMs, battery and Auto Part Manufacturers; end of life of EV vehicle - sustainability and circular economy. "

# Row 6
$ws.Range("A6").Value = "Bilateral consultations with industry members"
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = "2021-Second Half"
$ws.Range("D6").Value = "Convenings/workshops"
$ws.Range("E6").Value = "Unlocking policy and regulatory barriers"
$ws.Range("F6").Value = "First kick off platform meeting was held on 10 August 2021. Bilateral consultations with industry members were held between April and July 2021 to identify barriers faced by the different players across EV value chain and potential solutions required to address these challenges. 2nd meeting of the Platform discusses findings of the horizon scanning was held on August 10, 2021. Horizon screening comprised of barriers & solutions for scaling EVs in India (assessing member perspectives sought during consultations supported with a detailed review of literature)"

# Row 7
$ws.Range("A7").Value = "Bilateral consultations with industry members"
$ws.Range("B7").Value = "Completed"
$ws.Range("C7").Value = "2022-First Half"
$ws.Range("D7").Value = "Convenings/workshops"
$ws.Range("E7").Value = "Unlocking policy and regulatory barriers"
$ws.Range("F7").Value = "CII had had about 40 bilateral Consultations with industry members for consensus building around emerging solutions to address challenges faced by EV players are ongoing with partners. They have identified solutions with goal, scope, timelines and execution agency; members' inputs and Government feedback to finalise. "

# 5. Alignment tweaks
#    Title row (row1, style used by A1): vertical alignment top -> center
$ws.Range("A1").VerticalAlignment = -4108
#    Outputs row (row3): horizontal alignment now centered too
$ws.Range("A3:F3").HorizontalAlignment = -4108

Write-Output "edit complete"
